# ISIS-2062: removes testing.adoc as a category of domain services
# - moves SudoService to application-layer-api
# Also: isis-core-applib -> isis-applib.

$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Slide 1: shrink the "Called by framework" rectangle now that the
#    "Testing" / "Bootstrapping SPI" boxes it used to describe are gone.
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)

$calledByFramework = $s.Shapes.Item("Rectangle 21")
$calledByFramework.Height = 1800200 / 12700

# ------------------------------------------------------------------
# 2. Slide 1: remove the "Testing" category box and the
#    "Bootstrapping SPI" category box.
# ------------------------------------------------------------------
$s.Shapes.Item("Rectangle 19").Delete()
$s.Shapes.Item("Rectangle 16").Delete()

# ------------------------------------------------------------------
# 3. Refresh the cached "datetimeFigureOut" text on every slide layout
#    and on the slide master.
# ------------------------------------------------------------------
$newDate = "20/02/2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
